# Time and Payment.xlsx -- extend the D (Hours) shared formula down to row 49
# and scroll/select to reflect the author's new working area.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Extend the "C-B" hours formula from D40 through D49 -------------------
# D2:D8, D10:D24 and D32:D38 already hold "=Cn-Bn"; D39 is a text note and
# D40 (the last data row, A40/B40/C40 already filled in) is still missing
# its Hours formula. Continue the pattern down through D49 (new blank rows).
$ws.Range("D40:D49").Formula = "=C40-B40"

# Match the source formatting used by the rest of the D column (time format)
# and the sheet's standard row height for the freshly added rows 41:49.
$ws.Range("D40:D49").NumberFormat = "h:mm"
for ($r = 41; $r -le 49; $r++) {
    $ws.Rows.Item($r).RowHeight = 15.75
}

# --- Move the viewport / selection to where the new entries were made ------
$excel.ActiveWindow.ScrollRow = 24
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E47").Select() | Out-Null
